$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gains two
#    trailing spaces, followed by a red (C00000) parenthetical note
#    appended as three separate runs:
#      "(This is a change – Ve" / "rsion for branch alternate" / ")"
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
# Exclude the paragraph mark so we only touch the visible text.
$r1.End = $r1.End - 1
$r1.InsertAfter("  ")

$dash = [char]0x2013

$rNote1 = $d.Range($r1.End, $r1.End)
$rNote1.InsertAfter("(This is a change " + $dash + " Ve")
$rNote1.Font.Color = 192

$rNote2 = $d.Range($rNote1.End, $rNote1.End)
$rNote2.InsertAfter("rsion for branch alternate")
$rNote2.Font.Color = 192

$rNote3 = $d.Range($rNote2.End, $rNote2.End)
$rNote3.InsertAfter(")")
$rNote3.Font.Color = 192

# ------------------------------------------------------------------
# 2) Append a new, otherwise empty paragraph after the document's
#    final paragraph, shaded with fill F9F9F9 (clear/auto pattern).
#    InsertXML lets us drop in the exact WordprocessingML for the
#    paragraph without inheriting the preceding run's direct
#    character formatting.
# ------------------------------------------------------------------
$endPos = $d.Content.End
$rEnd = $d.Range($endPos, $endPos)
$rEnd.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>')
